$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster changes from "ECs" to "FAPs"; edge-weight columns recomputed with new TPM values
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8326193333333333
$ws.Range("N2").Value = 2.497858
$ws.Range("O2").Value = 0.9388124812781204
$ws.Range("P2").Value = 0.9388124812781203
$ws.Range("Q2").Value = 0.1420703919262222
$ws.Range("R2").Value = 1.278633527336
$ws.Range("S2").Value = 0.9388124812781204
$ws.Range("T2").Value = 0.9388124812781203

# Row 3: Target cluster changes from "FAPs" to "MuSCs"; edge-weight columns recomputed with new TPM values
$ws.Range("D3").Value = "MuSCs"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05426633333333333
$ws.Range("N3").Value = 0.162799
$ws.Range("O3").Value = 0.0611875187218796
$ws.Range("P3").Value = 0.06118751872187959
$ws.Range("Q3").Value = 0.009259500634222224
$ws.Range("R3").Value = 0.083335505708
$ws.Range("S3").Value = 0.0611875187218796
$ws.Range("T3").Value = 0.06118751872187959

# Row 4 (old "MuSCs" row) no longer exists in the updated dataset
$ws.Rows(4).Delete()
